$d = $word.ActiveDocument

# 1. Update the "Association Definition" text: the customer/order
#    relationship line is reworded.
$d.Content.Find.Execute("A customer orders many menu items.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A customer makes orders.", 2)

# 2. Append the remaining association-description paragraphs after it, using
#    unique placeholder tokens for what will end up being blank paragraphs
#    (so the paragraph break machinery doesn't leave a stray empty run
#    behind once the placeholder text is stripped out again below).
$cr = [char]13
$parts = @(
  "@@P1@@",
  "Each order has one to many menu items on it.",
  "@@P2@@",
  "@@P3@@",
  "A customer uses his/her miming’s account to pay for his orders when he orders to go or after he is finished eating.",
  "@@P4@@",
  "A staff works in a shift.",
  "@@P5@@",
  "@@P6@@",
  "@@P7@@"
)
$text = $cr + ($parts -join $cr)

$last = $d.Paragraphs($d.Paragraphs.Count)
$r = $last.Range
$r.Collapse(0)
$r.InsertAfter($text)

foreach ($ph in @("@@P1@@", "@@P2@@", "@@P3@@", "@@P4@@", "@@P5@@", "@@P6@@", "@@P7@@")) {
    $d.Content.Find.Execute($ph, $true, $false, $false, $false, $false,
                             $true, 1, $false, "", 2)
}

# 3. Relocate the "_GoBack" bookmark: it now sits mid-sentence in the Menu
#    Items description (splitting that run into "...description of t" /
#    "he food...") rather than at the end of the Association Definition
#    section.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$menuItemsPara = $d.Paragraphs(14)
$splitPoint = $menuItemsPara.Range.Start + 33
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange)
